$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) for the season record columns AD:AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new headers the same formatting as the rest of the header row
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Populate the season record for every data row (2-62)
for ($row = 2; $row -le 62; $row++) {
    $ws.Cells.Item($row, 30).Value = 68   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 93   # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
